# Generate Report for Handback
# Update the timestamp values recorded on the handback status report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file.
$wsOverview.Range("G2").Value = "2016-08-20 07:08:43"

# zh-cn sheet: handoff / handback datetimes for the first file.
$wsZhCn.Range("H2").Value = "2016-08-20 07:08:38"
$wsZhCn.Range("K2").Value = "2016-08-20 07:08:56"

# de-de sheet: "Latest HO Xliff Generate Date" (shared with Overview) and
# the handback datetime for the first file.
$wsDeDe.Range("H2").Value = "2016-08-20 07:08:43"
$wsDeDe.Range("K2").Value = "2016-08-20 07:09:07"
